# "clear content template excel"
#
# The "Pengguna" import template ships five sample/demo data rows
# (rows 9-13) that were filled in with example values (name "Bayu",
# password "pass", email "pass@mail.com", a phone number and some
# sample level codes) plus live mailto: hyperlinks on column E.
#
# This clears that sample data back out so the template only keeps its
# headers/instructions: row 9 keeps column E's formatting (it's the
# "Kode Level" example cell) but empty, columns A-D and F on every
# data row (9-13) are wiped completely (value + formatting), and the
# mailto hyperlinks that decorated column E are removed. The cursor
# is also left parked at K3 instead of C14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlinks that lived on E9:E13.
$ws.Hyperlinks.Delete()

# Columns, by index, for the A-D / E / F groups used below.
$colA = 1
$colD = 4
$colE = 5
$colF = 6

for ($r = 9; $r -le 13; $r++) {
    # A:D fully cleared (value + style).
    $abcd = $ws.Range($ws.Cells.Item($r, $colA), $ws.Cells.Item($r, $colD))
    $abcd.Clear()

    # F fully cleared (value + style) too.
    $ws.Cells.Item($r, $colF).Clear()

    # E keeps its style (part of the "Kode Level" column formatting),
    # only the sample value is removed.
    $ws.Cells.Item($r, $colE).ClearContents()
}

# Move the active selection from C14 to K3.
$ws.Range("K3").Select()
